$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of Channels input changes 105 -> 90 (drives the whole recompute cascade)
$ws.Range("B7").Value = 90

# New "Other Analyses" mini table in columns U:W
$ws.Range("U1").WrapText = $true
$ws.Range("U1").Value = "Top of Channel radius as a function of chamber radius"

$ws.Range("U2").WrapText = $true
$ws.Range("U2").Value = "Chamber radius"
$ws.Range("V2").Value = 78.7
$ws.Range("W2").Value = "mm"

$ws.Range("U4").WrapText = $true
$ws.Range("U4").Value = "Channel top radius"
$ws.Range("V4").Formula = "=(V2/1000+B5+V3)*1000"
$ws.Range("W4").Value = "mm"

$ws.Range("U3").WrapText = $true
$ws.Range("U3").Value = "channel height"
$ws.Range("V3").NumberFormat = "0.00E+00"
$ws.Range("V3").Formula = "=(((2*PI()*V2/1000)-B6*B7)/B7)"
$ws.Range("W3").Value = "m"

# Column widths for the new columns
$ws.Columns.Item(9).ColumnWidth = 11.28515625
$ws.Columns.Item(21).ColumnWidth = 22.42578125
$ws.Range("U1:U4").WrapText = $true

# Row heights: header row grows (wraps the long title), row 2 shrinks back
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(2).RowHeight = 30

# Selection moves to V3 as last edited cell
$ws.Range("V3").Select()
